# [WIP] Cambio de forma de PCB y componentes de Mouser
#
# The PCB shape / library components changed, so the BoM now sources three
# of the parts from Digi-Key instead of Farnell, with new supplier part
# numbers. The "Report Date" time stamp was also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report time stamp (row 6, "Report Date:") updated. Leading apostrophe
# keeps it a literal text label (matches the original "15:18" text cell)
# instead of letting Excel reinterpret "12:11" as a time value.
$ws.Range("D6").Value = "'12:11"

# Supplier 1 changed from Farnell to Digi-Key for the three sourced parts
# (rows 9, 11 and 12).
$ws.Range("F9").Value = "Digi-Key"
$ws.Range("F11").Value = "Digi-Key"
$ws.Range("F12").Value = "Digi-Key"

# New Digi-Key supplier part numbers replace the old Farnell numeric codes.
$ws.Range("G9").Value = "QLS6B-FKW-CNSNSF043CT-ND"
$ws.Range("G11").Value = "P16063CT-ND"
$ws.Range("G12").Value = "401-1910-1-ND"

# Leave the active selection on G12, matching the last-edited cell.
$ws.Range("G12").Select()
